$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new column K header (new question column)
$ws.Range("K1").Value = "Does it differ from PostgresSQL STIG tailored for CMS ARS 3.1?"

# Fill column K with "NO" for every data row (rows 2-81)
for ($r = 2; $r -le 81; $r++) {
    $ws.Range("K$r").Value = "NO"
}

# Header row grows taller to accommodate the wrapped new header text
$ws.Rows.Item(1).RowHeight = 34

# Update the active selection to I2, matching the saved view state
$null = $ws.Range("I2").Select()
